# ddCT2 model edit: add rate (E) and % of optimum rate (F) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Write the two new header strings FIRST, in the order they need to
# land in the shared-string table (E1's text before F1's text).
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "rate = ddCT2 of treat - noninjected / time; time = 4hr"
$ws.Range("F1").Value = "% of optimum rate"

# ---------------------------------------------------------------------
# Column F -- "% of optimum rate". Bold+red cells (header + the three
# "Topt" rows) are styled first so the bold-red font is registered
# before the plain-red font, matching the font/style table ordering
# produced by the original authoring session.
# ---------------------------------------------------------------------
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Font.Color = 255

$ws.Range("F2").Font.Bold = $true
$ws.Range("F2").Font.Color = 255

$ws.Range("F6").Formula = "=100"
$ws.Range("F6").Font.Bold = $true
$ws.Range("F6").Font.Color = 255

$ws.Range("F10").Value = 100
$ws.Range("F10").Font.Bold = $true
$ws.Range("F10").Font.Color = 255

# Plain (non-bold) red cells.
$ws.Range("F7").Formula = "=E7/E6*100"
$ws.Range("F7").Font.Color = 255

$ws.Range("F8").Formula = "=E8/E6*100"
$ws.Range("F8").Font.Color = 255

$ws.Range("F9").Formula = "=E9/E6*100"
$ws.Range("F9").Font.Color = 255

$ws.Range("F11").Formula = "=E11/E10*100"
$ws.Range("F11").Font.Color = 255

$ws.Range("F12").Formula = "=E12/E10*100"
$ws.Range("F12").Font.Color = 255

$ws.Range("F13").Formula = "=E13/E10*100"
$ws.Range("F13").Font.Color = 255

# ---------------------------------------------------------------------
# Column E -- "rate = ddCT2 of treat - noninjected / time; time = 4hr".
# The header and the three "Topt" rows reuse the workbook's existing
# plain-bold style; the rest stay on the default (unstyled) cell style.
# Resetting to "Normal" after writing the formula avoids inheriting
# column D's numeric format, which Excel otherwise applies to any new
# formula cell next to a formatted column.
# ---------------------------------------------------------------------
$ws.Range("E1").Style = "Normal"
$ws.Range("E1").Font.Bold = $true

$ws.Range("E6").Formula = "=C6-C2/4"
$ws.Range("E6").Style = "Normal"
$ws.Range("E6").Font.Bold = $true

$ws.Range("E10").Formula = "=C10-C2/4"
$ws.Range("E10").Style = "Normal"
$ws.Range("E10").Font.Bold = $true

$ws.Range("E7").Formula = "=C7-C3/4"
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Formula = "=C8-C4/4"
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Formula = "=C9-C5/4"
$ws.Range("E9").Style = "Normal"

$ws.Range("E11").Formula = "=C11-C3/4"
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Formula = "=C12-C4/4"
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Formula = "=C13-C5/4"
$ws.Range("E13").Style = "Normal"

# --- Column F width (matches new <col> entry in the diff) ---
$ws.Columns("F").ColumnWidth = 8.7265625

# --- Selection moves to F14, mirroring the source file's saved state ---
$ws.Range("F14").Select()
